# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates (and clears) per the target diff,
# row by row, sheet by sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 8
$ws.Range("H8").Value = 544.2
$ws.Range("I8").Value = 40.333332
$ws.Range("J8").Value = 1300
$ws.Range("K8").Value = 120.999996
$ws.Range("L8").Value = 3900
$ws.Range("M8").Value = 18.000004
$ws.Range("N8").Value = -4178

# ALC row 17
$ws.Range("H17").Value = 3731688.2
$ws.Range("J17").Value = 3820525
$ws.Range("L17").Value = 11461575
$ws.Range("N17").Value = -11461911

# ALC row 28
$ws.Range("H28").Value = 176.57143
$ws.Range("I28").Value = 178.61539
$ws.Range("J28").Value = 150
$ws.Range("K28").Value = 178.61539
$ws.Range("L28").Value = 150
$ws.Range("M28").Value = 306.38461
$ws.Range("N28").Value = -1120

# ALC row 99
$ws.Range("H99").Value = 1057.9231
$ws.Range("I99").Value = 812.75
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2438.25
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = -940.25
$ws.Range("N99").Value = -14996

# ALC row 113
$ws.Range("H113").Value = 6669346.5
$ws.Range("I113").Value = 8335945.5
$ws.Range("J113").Value = 2950
$ws.Range("K113").Value = 8335945.5
$ws.Range("L113").Value = 2950
$ws.Range("M113").Value = -8332691.5
$ws.Range("N113").Value = -9458

# ALC row 123
$ws.Range("H123").Value = 29598.13
$ws.Range("J123").Value = 29598.13
$ws.Range("L123").Value = 29598.13
$ws.Range("N123").Value = -39398.13

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 7396.9023
$ws.Range("I32").Value = 5293.0835
$ws.Range("K32").Value = 5293.0835
$ws.Range("M32").Value = -5006.0835

# ARM row 63
$ws.Range("H63").Value = 142860620
$ws.Range("I63").Value = 142860620
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 142860620
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -142859934
$ws.Range("N63").ClearContents()

# ARM row 66
$ws.Range("H66").Value = 142860620
$ws.Range("I66").Value = 142860620
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 714303100
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -714299668
$ws.Range("N66").ClearContents()

# ARM row 74
$ws.Range("H74").Value = 1270.9231
$ws.Range("I74").Value = 1147.3077
$ws.Range("J74").Value = 1518.1538
$ws.Range("K74").Value = 1147.3077
$ws.Range("L74").Value = 1518.1538
$ws.Range("M74").Value = -273.3077000000001
$ws.Range("N74").Value = -3266.1538

# ARM row 77
$ws.Range("H77").Value = 1270.9231
$ws.Range("I77").Value = 1147.3077
$ws.Range("J77").Value = 1518.1538
$ws.Range("K77").Value = 5736.538500000001
$ws.Range("L77").Value = 7590.769
$ws.Range("M77").Value = -1368.538500000001
$ws.Range("N77").Value = -16326.769

$ws = $wb.Worksheets.Item("BSM")
# BSM row 110
$ws.Range("H110").Value = 43000
$ws.Range("J110").Value = 43000
$ws.Range("L110").Value = 43000
$ws.Range("N110").Value = -51180

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 3425.4253
$ws.Range("I31").Value = 1504.174
$ws.Range("J31").Value = 5580.9756
$ws.Range("K31").Value = 1504.174
$ws.Range("L31").Value = 5580.9756
$ws.Range("M31").Value = -1209.174
$ws.Range("N31").Value = -6170.9756

# CRP row 34
$ws.Range("H34").Value = 3425.4253
$ws.Range("I34").Value = 1504.174
$ws.Range("J34").Value = 5580.9756
$ws.Range("K34").Value = 1504.174
$ws.Range("L34").Value = 5580.9756
$ws.Range("M34").Value = -1302.174
$ws.Range("N34").Value = -5984.9756

# CRP row 100
$ws.Range("H100").Value = 52780
$ws.Range("J100").Value = 52780
$ws.Range("L100").Value = 52780
$ws.Range("N100").Value = -54944

# CRP row 132
$ws.Range("H132").Value = 3302.5
$ws.Range("I132").Value = 3120.6667
$ws.Range("J132").Value = 3484.3333
$ws.Range("K132").Value = 9362.000100000001
$ws.Range("L132").Value = 10452.9999
$ws.Range("M132").Value = -6832.000100000001
$ws.Range("N132").Value = -15512.9999

# CRP row 134
$ws.Range("H134").Value = 3981.5833
$ws.Range("I134").Value = 4237.684
$ws.Range("J134").Value = 3008.4
$ws.Range("K134").Value = 12713.052
$ws.Range("L134").Value = 9025.200000000001
$ws.Range("M134").Value = -10178.052
$ws.Range("N134").Value = -14095.2

$ws = $wb.Worksheets.Item("CUL")
# CUL row 136
$ws.Range("H136").Value = 4320.566
$ws.Range("I136").Value = 11921.111
$ws.Range("J136").Value = 2765.9092
$ws.Range("K136").Value = 35763.333
$ws.Range("L136").Value = 8297.7276
$ws.Range("M136").Value = -30663.333
$ws.Range("N136").Value = -18497.7276

$ws = $wb.Worksheets.Item("GSM")
# GSM row 40
$ws.Range("H40").Value = 12800
$ws.Range("J40").Value = 12800
$ws.Range("L40").Value = 12800
$ws.Range("N40").Value = -13102

# GSM row 43
$ws.Range("H43").Value = 1600
$ws.Range("I43").Value = 1600
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1600
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1449
$ws.Range("N43").ClearContents()

# GSM row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# GSM row 102
$ws.Range("H102").Value = 943187.7
$ws.Range("I102").Value = 1884216.6
$ws.Range("K102").Value = 1884216.6
$ws.Range("M102").Value = -1882594.6

# GSM row 132
$ws.Range("H132").Value = 2435
$ws.Range("I132").Value = 1754.5385
$ws.Range("J132").Value = 3540.75
$ws.Range("K132").Value = 5263.6155
$ws.Range("L132").Value = 10622.25
$ws.Range("M132").Value = -2733.6155
$ws.Range("N132").Value = -15682.25

$ws = $wb.Worksheets.Item("LTW")
# LTW row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# LTW row 61
$ws.Range("H61").Value = 1743.1111
$ws.Range("I61").Value = 1848
$ws.Range("J61").Value = 1533.3334
$ws.Range("K61").Value = 1848
$ws.Range("L61").Value = 1533.3334
$ws.Range("M61").Value = -1646
$ws.Range("N61").Value = -1937.3334

# LTW row 113
$ws.Range("H113").Value = 1743.1111
$ws.Range("I113").Value = 1848
$ws.Range("J113").Value = 1533.3334
$ws.Range("K113").Value = 1848
$ws.Range("L113").Value = 1533.3334
$ws.Range("M113").Value = 322
$ws.Range("N113").Value = -5873.3334

# LTW row 122
$ws.Range("H122").Value = 3706462.2
$ws.Range("I122").Value = 4207362.5
$ws.Range("J122").Value = 2003401
$ws.Range("K122").Value = 12622087.5
$ws.Range("L122").Value = 6010203
$ws.Range("M122").Value = -12619637.5
$ws.Range("N122").Value = -6015103

# LTW row 132
$ws.Range("H132").Value = 15154486
$ws.Range("I132").Value = 25643138
$ws.Range("J132").Value = 4210.778
$ws.Range("K132").Value = 76929414
$ws.Range("L132").Value = 12632.334
$ws.Range("M132").Value = -76926884
$ws.Range("N132").Value = -17692.334

$ws = $wb.Worksheets.Item("WVR")
# WVR row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = 0

# WVR row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = 0

# WVR row 132
$ws.Range("H132").Value = 1673.45
$ws.Range("I132").Value = 1325.7391
$ws.Range("J132").Value = 2143.8823
$ws.Range("K132").Value = 3977.2173
$ws.Range("L132").Value = 6431.646900000001
$ws.Range("M132").Value = -1447.2173
$ws.Range("N132").Value = -11491.6469

# WVR row 136
$ws.Range("H136").Value = 1881.4916
$ws.Range("I136").Value = 1776.027
$ws.Range("J136").Value = 2058.8635
$ws.Range("K136").Value = 5328.081
$ws.Range("L136").Value = 6176.5905
$ws.Range("M136").Value = -2778.081
$ws.Range("N136").Value = -11276.5905

